# Edit script for "Квест с паучьей королевой.docx"
# Applies spell/grammar-checker proofErr markup splits and appends
# "(Сделано)"/"(Исправлено)" status markers, plus two new list items
# at the end of the document.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($paraIndex, $innerXml, $pPrXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $xml = "<w:p $wNs>" + $pPrXml + $innerXml + "</w:p>"
    $null = $p.Range.InsertXML($xml)
}

$pPr1 = '<w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>'
$pPr1 = '<w:pPr>' + $pPr1 + '</w:pPr>'
$pPr2 = '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$pPr3 = '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$pPr5 = '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>'
$pPrNone = ''

# --- Paragraph 1: intro sentence -----------------------------------------
$inner = '<w:r><w:t xml:space="preserve">Для реализации этого </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квеста</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> мне нужно пройтись по следующим пунктам:</w:t></w:r>'
Set-ParaXml 1 $inner $pPrNone

# --- Paragraph 2: "Квесты и истории..." -----------------------------------
$inner = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Квесты</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> и истории, длящиеся на протяжении всей игры</w:t></w:r>' +
    '<w:r><w:t>(Сделано)</w:t></w:r>'
Set-ParaXml 2 $inner $pPr1

# --- Paragraph 5: "Анимации во время диалогов..." -------------------------
$inner = '<w:r><w:t>Анимации во время диалогов</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, расширенный редактор диалогов, который включает в себя управление камерой, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>анимациями</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> персонажей, сменой места действия и возможность ответа на </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>вопрос.</w:t></w:r>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>Сделано)</w:t></w:r>'
Set-ParaXml 5 $inner $pPr1

# --- Paragraph 8: "Первая часть квеста" ------------------------------------
$inner = '<w:r><w:t xml:space="preserve">Первая часть </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квеста</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>(Сделано)</w:t></w:r>'
Set-ParaXml 8 $inner $pPr1

# --- Paragraph 9: "Вторая часть квеста" ------------------------------------
$inner = '<w:r><w:t xml:space="preserve">Вторая часть </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квеста</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> (Сделано)</w:t></w:r>'
Set-ParaXml 9 $inner $pPr1

# --- Paragraph 10: "Третья часть квеста" -----------------------------------
$inner = '<w:r><w:t xml:space="preserve">Третья часть </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квеста</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> (Сделано)</w:t></w:r>'
Set-ParaXml 10 $inner $pPr1

# --- Paragraph 12: "Что будет, если подойти к пауку лазутчику..." ---------
$inner = '<w:r><w:t>Что будет, если подойти к пауку лазутчику (вторая</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> часть </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квеста</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>) с другой стороны?</w:t></w:r>'
Set-ParaXml 12 $inner $pPr2

# --- Paragraph 14: "Странное сохранение игры..." ---------------------------
$inner = '<w:r><w:t xml:space="preserve">Странное сохранение игры после выполнения </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квеста</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, дух не следует за героем, да и почему-то дублируются игровые </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>эффекты.</w:t></w:r>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>Исправл</w:t></w:r>' +
    '<w:r><w:t>ено)</w:t></w:r>'
Set-ParaXml 14 $inner $pPr3

# --- Paragraph 15: "Нехорошо, когда умершие..." ----------------------------
$inner = '<w:r><w:t xml:space="preserve">Нехорошо, когда умершие </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квестовые</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> персонажи могут стать союзниками и нужно ждать их смерти, чтобы запустился </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>тригер</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квестового</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> события</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:r><w:t>Исправлено</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>'
Set-ParaXml 15 $inner $pPr3

# --- Paragraph 16: "Телепортация паука-героя и его развороты" -------------
$inner = '<w:r><w:t>Телепортация паука-героя и его развороты</w:t></w:r>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:r><w:t>Исправлено</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>'
Set-ParaXml 16 $inner $pPr3

# --- Paragraph 17: "Паук-герой немного неудобен..." ------------------------
$inner = '<w:r><w:t>Паук-герой немного неудобен, когда подходит к краям – он с некоторой вероятностью падает</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:r><w:t>Исправлено</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>'
Set-ParaXml 17 $inner $pPr3

# --- Paragraph 18: "Странная загрузка игры..." -----------------------------
$inner = '<w:r><w:t>Странная загрузка игры после превращения в паука и смерти в этом обличии</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:r><w:t>Исправлено</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>'
Set-ParaXml 18 $inner $pPr3

# --- Paragraph 19: "Исправить внешний вид героев" --------------------------
$inner = '<w:r><w:t>Исправить внешний вид героев</w:t></w:r>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:r><w:t>Сделано)</w:t></w:r>'
Set-ParaXml 19 $inner $pPr3

# --- Paragraph 20: "Если герой дошёл до скрытого убежища..." --------------
$inner = '<w:r><w:t>Если герой дошёл до скрытого убежища, то паук вор автоматически заходит в неё – не ждёт игрока</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (Сделано)</w:t></w:r>'
Set-ParaXml 20 $inner $pPr3

# --- Paragraph 21: "Анимация ухода паука..." -------------------------------
$inner = '<w:r><w:t>Анимация ухода паука после выполнения задания – её нет.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (Исправлено)</w:t></w:r>'
Set-ParaXml 21 $inner $pPr3

# --- Append two new list items after the last paragraph (23) --------------
$p23 = $d.Paragraphs.Item(23)
$insertPoint = $d.Range($p23.Range.End, $p23.Range.End)

$newPara1Inner = '<w:r><w:t xml:space="preserve">В первой части </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квеста</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, если сразу прийти к пауку-вору и уйти, а потом взять </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>квест</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, то с пауком вором происходит что-то странное</w:t></w:r>'

$newPara2Inner = '<w:r><w:t xml:space="preserve">Исправить </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>тайминги</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> диалогов, дух должен находиться рядом с героем во время перехода в другую часть локации.</w:t></w:r>'

$newParasXml = ('<w:p ' + $wNs + '>' + $pPr5 + $newPara1Inner + '</w:p>') +
    ('<w:p ' + $wNs + '>' + $pPr5 + $newPara2Inner + '</w:p>')

$null = $insertPoint.InsertXML($newParasXml)
